$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Taul1")

# Copy formatting from the row above (row 17) into row 18 so the new
# entry matches the rest of the table (borders, number formats, etc.)
$ws.Range("A17:F17").Copy()
$ws.Range("A18:F18").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the new tracking entry for row 18
$ws.Range("A18").Value = 8.12
$ws.Range("B18").Value = 0.375
$ws.Range("C18").Value = 0.625
$ws.Range("E18").Value = "6hr"
$ws.Range("F18").Value = "After setting up basic layout then try to insert function in every pages until so far function 2 is inserted"
